$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.777.66'
$ws.Range("E2").Value = '  +0.40%  '

$ws.Range("D3").Value = '1.751.05'

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = "'236.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.65%  '

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = "'0.5064"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'40.42"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.2615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.51%  '

$ws.Range("D10").Value = "'0.06199"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = '1.752.46'
$ws.Range("E11").Value = '  +0.78%  '

$ws.Range("D12").Value = "'0.06935"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.51%  '

$ws.Range("E13").Value = '  +5.58%  '

$ws.Range("D14").Value = "'0.6059"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.35%  '

$ws.Range("D15").Value = "'78.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.12%  '

$ws.Range("D16").Value = "'4.461"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.67%  '

$ws.Range("D17").Value = "'0.9998"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.22%  '

$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("D19").Value = '25.803.80'
$ws.Range("E19").Value = '  +0.39%  '

$ws.Range("D20").Value = "'11.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.03%  '

$ws.Range("D21").Value = "'0.000006725"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.17%  '

$ws.Range("D22").Value = '1.974.91'
$ws.Range("E22").Value = '  +0.91%  '

$ws.Range("E23").Value = '  +3.11%  '

$ws.Range("D24").Value = "'8.192"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.67%  '

$ws.Range("E25").Value = '  +1.24%  '

$ws.Range("D26").Value = "'137.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.52%  '

$ws.Range("D27").Value = "'1.456"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").Value = "'15.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.99%  '

$ws.Range("D29").Value = "'1.807"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.23%  '

$ws.Range("D30").Value = "'102.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.30%  '

$ws.Range("D31").Value = "'0.08264"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.31%  '

$ws.Range("E32").Value = '  +0.14%  '

$ws.Range("D33").Value = "'3.397"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.46%  '

$ws.Range("D34").Value = "'0.04348"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("D35").Value = "'0.9992"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("D36").Value = "'2.646"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.30%  '

$ws.Range("D37").Value = "'0.9999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.83%  '

$ws.Range("D38").Value = "'0.6013"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.75%  '

$ws.Range("D39").Value = "'2.704"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.23%  '

$ws.Range("D40").Value = "'1.961"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.65%  '

$ws.Range("D41").Value = "'0.01548"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.51%  '

$ws.Range("D42").Value = "'1.000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("D43").Value = "'103.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.59%  '

$ws.Range("D44").Value = "'0.7553"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.53%  '

$ws.Range("D45").Value = "'0.3813"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.76%  '

$ws.Range("D46").Value = "'4.873"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.61%  '

$ws.Range("D47").Value = "'0.05484"
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = '  +2.67%  '

$ws.Range("D49").Value = "'5.937"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.87%  '

$ws.Range("D50").Value = "'30.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.11%  '

$ws.Range("E51").Value = '  +0.48%  '
